$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row: "_old" -> "_FV2310", "_new" -> "_FV2404" ---
$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le 21; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $text = $cell.Value()
    if ($text -ne $null) {
        $text = $text -replace "_old$", "_FV2310"
        $text = $text -replace "_new$", "_FV2404"
        $cell.Value = $text
    }
}

# --- Turn the data range into an Excel Table (ListObject) with autofilter ---
$dataRange = $ws.Range("A1:U51")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# --- Freeze the header row (split/freeze pane below row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
